$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.070935
$ws.Cells.Item(2, 8).Value = 6.212805
$ws.Cells.Item(2, 9).Value = 0.09632749399019591
$ws.Cells.Item(2, 10).Value = 0.09632749399019594
$ws.Cells.Item(2, 13).Value = 22.91402233333334
$ws.Cells.Item(2, 14).Value = 68.74206700000001
$ws.Cells.Item(2, 15).Value = 0.191813973987922
$ws.Cells.Item(2, 16).Value = 0.191813973987922
$ws.Cells.Item(2, 17).Value = 47.45345084088167
$ws.Cells.Item(2, 18).Value = 427.081057567935
$ws.Cells.Item(2, 19).Value = 0.01847695942655715
$ws.Cells.Item(2, 20).Value = 0.01847695942655716
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.070935
$ws.Cells.Item(3, 8).Value = 6.212805
$ws.Cells.Item(3, 9).Value = 0.09632749399019591
$ws.Cells.Item(3, 10).Value = 0.09632749399019594
$ws.Cells.Item(3, 15).Value = 0.07776078244711707
$ws.Cells.Item(3, 16).Value = 0.07776078244711705
$ws.Cells.Item(3, 17).Value = 19.23747988994333
$ws.Cells.Item(3, 18).Value = 173.13731900949
$ws.Cells.Item(3, 19).Value = 0.007490501303847601
$ws.Cells.Item(3, 20).Value = 0.007490501303847602
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.070935
$ws.Cells.Item(4, 8).Value = 6.212805
$ws.Cells.Item(4, 9).Value = 0.09632749399019591
$ws.Cells.Item(4, 10).Value = 0.09632749399019594
$ws.Cells.Item(4, 13).Value = 6.020714333333333
$ws.Cells.Item(4, 14).Value = 18.062143
$ws.Cells.Item(4, 15).Value = 0.05039958178109668
$ws.Cells.Item(4, 16).Value = 0.05039958178109668
$ws.Cells.Item(4, 17).Value = 12.46850803790167
$ws.Cells.Item(4, 18).Value = 112.216572341115
$ws.Cells.Item(4, 19).Value = 0.004854865411126978
$ws.Cells.Item(4, 20).Value = 0.004854865411126979
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.070935
$ws.Cells.Item(5, 8).Value = 6.212805
$ws.Cells.Item(5, 9).Value = 0.09632749399019591
$ws.Cells.Item(5, 10).Value = 0.09632749399019594
$ws.Cells.Item(5, 13).Value = 81.23559966666666
$ws.Cells.Item(5, 14).Value = 243.706799
$ws.Cells.Item(5, 15).Value = 0.6800256617838641
$ws.Cells.Item(5, 16).Value = 0.6800256617838643
$ws.Cells.Item(5, 17).Value = 168.2336465956883
$ws.Cells.Item(5, 18).Value = 1514.102819361195
$ws.Cells.Item(5, 19).Value = 0.06550516784866417
$ws.Cells.Item(5, 20).Value = 0.0655051678486642
$ws.Cells.Item(6, 9).Value = 0.1492700517445119
$ws.Cells.Item(6, 10).Value = 0.1492700517445119
$ws.Cells.Item(6, 13).Value = 22.91402233333334
$ws.Cells.Item(6, 14).Value = 68.74206700000001
$ws.Cells.Item(6, 15).Value = 0.191813973987922
$ws.Cells.Item(6, 16).Value = 0.191813973987922
$ws.Cells.Item(6, 17).Value = 73.53434382083057
$ws.Cells.Item(6, 18).Value = 661.8090943874751
$ws.Cells.Item(6, 19).Value = 0.02863208182249757
$ws.Cells.Item(6, 20).Value = 0.02863208182249758
$ws.Cells.Item(7, 9).Value = 0.1492700517445119
$ws.Cells.Item(7, 10).Value = 0.1492700517445119
$ws.Cells.Item(7, 15).Value = 0.07776078244711707
$ws.Cells.Item(7, 16).Value = 0.07776078244711705
$ws.Cells.Item(7, 19).Value = 0.01160735601957489
$ws.Cells.Item(7, 20).Value = 0.01160735601957489
$ws.Cells.Item(8, 9).Value = 0.1492700517445119
$ws.Cells.Item(8, 10).Value = 0.1492700517445119
$ws.Cells.Item(8, 13).Value = 6.020714333333333
$ws.Cells.Item(8, 14).Value = 18.062143
$ws.Cells.Item(8, 15).Value = 0.05039958178109668
$ws.Cells.Item(8, 16).Value = 0.05039958178109668
$ws.Cells.Item(8, 17).Value = 19.32132523019722
$ws.Cells.Item(8, 18).Value = 173.891927071775
$ws.Cells.Item(8, 19).Value = 0.007523148180366058
$ws.Cells.Item(8, 20).Value = 0.00752314818036606
$ws.Cells.Item(9, 9).Value = 0.1492700517445119
$ws.Cells.Item(9, 10).Value = 0.1492700517445119
$ws.Cells.Item(9, 13).Value = 81.23559966666666
$ws.Cells.Item(9, 14).Value = 243.706799
$ws.Cells.Item(9, 15).Value = 0.6800256617838641
$ws.Cells.Item(9, 16).Value = 0.6800256617838643
$ws.Cells.Item(9, 17).Value = 260.6965477069528
$ws.Cells.Item(9, 18).Value = 2346.268929362575
$ws.Cells.Item(9, 19).Value = 0.1015074657220733
$ws.Cells.Item(9, 20).Value = 0.1015074657220734
$ws.Cells.Item(10, 7).Value = 1.134776333333333
$ws.Cells.Item(10, 8).Value = 3.404329
$ws.Cells.Item(10, 9).Value = 0.05278299919088877
$ws.Cells.Item(10, 10).Value = 0.05278299919088877
$ws.Cells.Item(10, 13).Value = 22.91402233333334
$ws.Cells.Item(10, 14).Value = 68.74206700000001
$ws.Cells.Item(10, 15).Value = 0.191813973987922
$ws.Cells.Item(10, 16).Value = 0.191813973987922
$ws.Cells.Item(10, 17).Value = 26.00229024533811
$ws.Cells.Item(10, 18).Value = 234.020612208043
$ws.Cells.Item(10, 19).Value = 0.01012451683380565
$ws.Cells.Item(10, 20).Value = 0.01012451683380565
$ws.Cells.Item(11, 7).Value = 1.134776333333333
$ws.Cells.Item(11, 8).Value = 3.404329
$ws.Cells.Item(11, 9).Value = 0.05278299919088877
$ws.Cells.Item(11, 10).Value = 0.05278299919088877
$ws.Cells.Item(11, 15).Value = 0.07776078244711707
$ws.Cells.Item(11, 16).Value = 0.07776078244711705
$ws.Cells.Item(11, 17).Value = 10.54124677601356
$ws.Cells.Item(11, 18).Value = 94.871220984122
$ws.Cells.Item(11, 19).Value = 0.004104447316989057
$ws.Cells.Item(11, 20).Value = 0.004104447316989057
$ws.Cells.Item(12, 7).Value = 1.134776333333333
$ws.Cells.Item(12, 8).Value = 3.404329
$ws.Cells.Item(12, 9).Value = 0.05278299919088877
$ws.Cells.Item(12, 10).Value = 0.05278299919088877
$ws.Cells.Item(12, 13).Value = 6.020714333333333
$ws.Cells.Item(12, 14).Value = 18.062143
$ws.Cells.Item(12, 15).Value = 0.05039958178109668
$ws.Cells.Item(12, 16).Value = 0.05039958178109668
$ws.Cells.Item(12, 17).Value = 6.832164135227445
$ws.Cells.Item(12, 18).Value = 61.489477217047
$ws.Cells.Item(12, 19).Value = 0.002660241084372758
$ws.Cells.Item(12, 20).Value = 0.002660241084372759
$ws.Cells.Item(13, 7).Value = 1.134776333333333
$ws.Cells.Item(13, 8).Value = 3.404329
$ws.Cells.Item(13, 9).Value = 0.05278299919088877
$ws.Cells.Item(13, 10).Value = 0.05278299919088877
$ws.Cells.Item(13, 13).Value = 81.23559966666666
$ws.Cells.Item(13, 14).Value = 243.706799
$ws.Cells.Item(13, 15).Value = 0.6800256617838641
$ws.Cells.Item(13, 16).Value = 0.6800256617838643
$ws.Cells.Item(13, 17).Value = 92.18423592587455
$ws.Cells.Item(13, 18).Value = 829.658123332871
$ws.Cells.Item(13, 19).Value = 0.0358937939557213
$ws.Cells.Item(13, 20).Value = 0.03589379395572131
$ws.Cells.Item(14, 7).Value = 15.08404533333333
$ws.Cells.Item(14, 8).Value = 45.252136
$ws.Cells.Item(14, 9).Value = 0.7016194550744034
$ws.Cells.Item(14, 10).Value = 0.7016194550744034
$ws.Cells.Item(14, 13).Value = 22.91402233333334
$ws.Cells.Item(14, 14).Value = 68.74206700000001
$ws.Cells.Item(14, 15).Value = 0.191813973987922
$ws.Cells.Item(14, 16).Value = 0.191813973987922
$ws.Cells.Item(14, 17).Value = 345.6361516450125
$ws.Cells.Item(14, 18).Value = 3110.725364805112
$ws.Cells.Item(14, 19).Value = 0.1345804159050616
$ws.Cells.Item(14, 20).Value = 0.1345804159050616
$ws.Cells.Item(15, 7).Value = 15.08404533333333
$ws.Cells.Item(15, 8).Value = 45.252136
$ws.Cells.Item(15, 9).Value = 0.7016194550744034
$ws.Cells.Item(15, 10).Value = 0.7016194550744034
$ws.Cells.Item(15, 15).Value = 0.07776078244711707
$ws.Cells.Item(15, 16).Value = 0.07776078244711705
$ws.Cells.Item(15, 17).Value = 140.1198100176942
$ws.Cells.Item(15, 18).Value = 1261.078290159248
$ws.Cells.Item(15, 19).Value = 0.05455847780670551
$ws.Cells.Item(15, 20).Value = 0.0545584778067055
$ws.Cells.Item(16, 7).Value = 15.08404533333333
$ws.Cells.Item(16, 8).Value = 45.252136
$ws.Cells.Item(16, 9).Value = 0.7016194550744034
$ws.Cells.Item(16, 10).Value = 0.7016194550744034
$ws.Cells.Item(16, 13).Value = 6.020714333333333
$ws.Cells.Item(16, 14).Value = 18.062143
$ws.Cells.Item(16, 15).Value = 0.05039958178109668
$ws.Cells.Item(16, 16).Value = 0.05039958178109668
$ws.Cells.Item(16, 17).Value = 90.81672794304978
$ws.Cells.Item(16, 18).Value = 817.350551487448
$ws.Cells.Item(16, 19).Value = 0.03536132710523088
$ws.Cells.Item(16, 20).Value = 0.03536132710523088
$ws.Cells.Item(17, 7).Value = 15.08404533333333
$ws.Cells.Item(17, 8).Value = 45.252136
$ws.Cells.Item(17, 9).Value = 0.7016194550744034
$ws.Cells.Item(17, 10).Value = 0.7016194550744034
$ws.Cells.Item(17, 13).Value = 81.23559966666666
$ws.Cells.Item(17, 14).Value = 243.706799
$ws.Cells.Item(17, 15).Value = 0.6800256617838641
$ws.Cells.Item(17, 16).Value = 0.6800256617838643
$ws.Cells.Item(17, 17).Value = 1225.361468052518
$ws.Cells.Item(17, 18).Value = 11028.25321247266
$ws.Cells.Item(17, 19).Value = 0.4771192342574053
$ws.Cells.Item(17, 20).Value = 0.4771192342574054
